# Apply cryptos.xlsx price/volume update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.103.04"
$ws.Range("E2").Value = "  -0.47%  "

$ws.Range("D3").Value = "1.661.91"
$ws.Range("E3").Value = "  -1.18%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.76"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5179"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -2.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2585"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -3.89%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06309"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.95"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -2.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07523"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -0.12%  "

$ws.Range("D12").Value = "1.646.74"
$ws.Range("E12").Value = "  -1.98%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.409"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -2.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5361"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -5.84%  "

$ws.Range("E15").Value = "  -0.78%  "

$ws.Range("D16").Value = "0.0₅7925"
$ws.Range("E16").Value = "  -3.09%  "

$ws.Range("D17").Value = "26.136.81"
$ws.Range("E17").Value = "  -0.39%  "

$ws.Range("E18").Value = "  -0.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.688"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -3.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "187.52"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.16"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -3.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.180"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.84%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.88"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +0.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1219"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -4.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.391"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -3.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.61"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -2.36%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.377"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +2.52%  "

$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06145"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -5.46%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.261"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -1.93%  "

$ws.Range("E31").Value = "  -2.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.390"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -2.98%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.629"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -1.83%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9849"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -2.65%  "

$ws.Range("E35").Value = "  -0.89%  "

$ws.Range("E36").Value = "  +1.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5864"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -4.22%  "

$ws.Range("D38").Value = "1.103.43"
$ws.Range("E38").Value = "  -0.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01590"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -2.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.973"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -3.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8439"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -2.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.83"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -0.27%  "

$ws.Range("D45").Value = "0.0₈107"
$ws.Range("E45").Value = "  -2.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "54.88"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -3.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9982"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -0.92%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05245"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -0.49%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.963"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4247"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.864"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -2.05%  "
